$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").Value = 0.9999999936462307
$ws.Range("E2").Value = 0.9999999936462307

# Row 3
$ws.Range("D3").Value = 0.0001584605617075297
$ws.Range("E3").Value = 0.0001584605617075297

# Row 4
$ws.Range("D4").Value = 0.9999999999999982
$ws.Range("E4").Value = 0.9999999999999982

# Row 5
$ws.Range("D5").Value = 0.999975568849983
$ws.Range("E5").Value = 0.999975568849983

# Row 6
$ws.Range("D6").Value = 0.9999999999999301
$ws.Range("E6").Value = 0.9999999999999301

# Row 7
$ws.Range("D7").Value = 0.9993481691473288
$ws.Range("E7").Value = 0.0006518308526711758

# Row 9
$ws.Range("D9").Value = 0.9016350470064394
$ws.Range("E9").Value = 0.09836495299356063

# Row 11
$ws.Range("D11").Value = 0.000001004439855873645
$ws.Range("E11").Value = 0.9999989955601442
$ws.Range("F11").Value = 10.76978397369385

$wb.Save()
